$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17; this shifts existing rows 17-69 down to 18-70
# and keeps formatting/styles consistent (mirrors Excel's native row insert).
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new weekly record.
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C17").Value = "Los Lagos"
$ws.Range("D17").Value = 44659
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100112031
$ws.Range("G17").Value = "Poroto verde"
$ws.Range("H17").Value = "Magnum"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 30000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 30000
$ws.Range("N17").Value = "$/saco 25 kilos"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 1200
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
